$wb = $excel.ActiveWorkbook

# --- Sheet "Informações Gerais" (sheet1): add 4 new occurrence rows ---
$wsGerais = $wb.Worksheets.Item("Informações Gerais")

# Insert 4 new rows above the current trailing blank row (row 38), pushing
# the blank placeholder row down to become the new trailing row (42).
# Excel's Insert inherits formatting from the row above, which matches the
# existing data rows exactly (same styles as row 37).
$wsGerais.Rows("38:41").Insert(-4121)  # xlShiftDown

$wsGerais.Range("A38").Value = 45415.40557880787
$wsGerais.Range("B38").Value = "0338.9/2024"
$wsGerais.Range("C38").Value = "Polícia Militar"
$wsGerais.Range("E38").Value = "SGT. BARRETTO"
$wsGerais.Range("F38").Value = 1104438
$wsGerais.Range("G38").Value = 13
$wsGerais.Range("H38").Value = "Apenas vítima(s)"

$wsGerais.Range("A39").Value = 45417.64596840278
$wsGerais.Range("B39").Value = "0376.9/2024"
$wsGerais.Range("C39").Value = "Polícia Militar"
$wsGerais.Range("E39").Value = "Sgt. P. Cavalcante"
$wsGerais.Range("F39").Value = 1104870
$wsGerais.Range("G39").Value = 25
$wsGerais.Range("H39").Value = "Apenas vítima(s)"

$wsGerais.Range("A40").Value = 45418.577854675925
$wsGerais.Range("B40").Value = "0382.9/2024"
$wsGerais.Range("C40").Value = "Polícia Militar"
$wsGerais.Range("E40").Value = "Sd. Talles"
$wsGerais.Range("F40").Value = 1252798
$wsGerais.Range("G40").Value = 25
$wsGerais.Range("H40").Value = "Apenas vítima(s)"

$wsGerais.Range("A41").Value = 45418.900479247684
$wsGerais.Range("B41").Value = "0385.9/2024"
$wsGerais.Range("C41").Value = "Polícia Militar"
$wsGerais.Range("E41").Value = "Cb. Damasceno"
$wsGerais.Range("F41").Value = 1182170
$wsGerais.Range("G41").Value = 6
$wsGerais.Range("H41").Value = "Apenas vítima(s)"

# --- Sheet "Vítimas" (sheet2): add the matching 4 victim summary rows ---
$wsVitimas = $wb.Worksheets.Item("Vítimas")

# Insert 4 rows before the current trailing blank row (row 39), shifting
# everything below down by 4 (rows 39..138 -> 43..142), inheriting the
# format of row 38 for the newly inserted rows.
$wsVitimas.Rows("39:42").Insert(-4121)  # xlShiftDown

$wsVitimas.Range("A39").Value = 45415.40603083333
$wsVitimas.Range("B39").Value = "0338.9/2024"
$wsVitimas.Range("C39").Value = 146204

$wsVitimas.Range("A40").Value = 45417.646333240744
$wsVitimas.Range("B40").Value = "0376.9/2024"
$wsVitimas.Range("C40").Value = 146212

$wsVitimas.Range("A41").Value = 45418.57809548611
$wsVitimas.Range("B41").Value = "0382.9/2024"
$wsVitimas.Range("C41").Value = 146425

$wsVitimas.Range("A42").Value = 45418.900763784724
$wsVitimas.Range("B42").Value = "0385.9/2024"
$wsVitimas.Range("C42").Value = 146422
